$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warning")

# --- Add the two new warning-message columns (B: "Warning Msg 2", C: "Warning Msg 3") ---
# Values are written column-by-column (B1,B2 then C1,C2) so the new shared-string
# table entries come out in the same order as the target workbook (98..101).

# Column B - header (row 1, bold/centered) then message (row 2, left/center + wrap).
$ws.Range("B1").Value = "Warning Msg 2"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4108

$ws.Range("B2").Value = "A Subject is typically considered a potential round trip if it is an operating company acquired either by a Private Equity firm or by a PE-owned operating company. The Buyer is not listed as a Private Equity/Hedge Fund/Family Office company or an Operating Company with Private Equity/Hedge Fund/Family Office ownership. If you still want to consider the Subject a round-trip candidate, no change is needed; otherwise, please change the selection."
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("B2").WrapText = $true

# Column C - header (row 1, bold/centered) then message (row 2, left/center + wrap).
$ws.Range("C1").Value = "Warning Msg 3"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4108

$ws.Range("C2").Value = "A Subject is typically considered a potential round trip if it is an operating company acquired either by a Private Equity firm or by a PE-owned operating company. The Subject is not listed as an Operating Company AND the Buyer is not listed as a Private Equity/Hedge Fund/Family Office company or an Operating Company with Private Equity/Hedge Fund/Family Office ownership. If you still want to consider the Subject a round-trip candidate no change is needed; otherwise, please change the selection."
$ws.Range("C2").HorizontalAlignment = -4131
$ws.Range("C2").VerticalAlignment = -4108
$ws.Range("C2").WrapText = $true

# Column widths for B and C to match column A's look-and-feel.
$ws.Columns.Item(2).ColumnWidth = 66
$ws.Columns.Item(3).ColumnWidth = 66

# Taller row 2 so the longer warning text still wraps cleanly.
$ws.Rows.Item(2).RowHeight = 100.8

# Make the "Warning" sheet the active tab/sheet and set the selection as seen after editing.
$ws.Activate()
$ws.Range("C6").Select() | Out-Null
